$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.013.78"
$ws.Range("E2").Value = "  +1.94%  "
$ws.Range("D3").Value = "1.647.64"
$ws.Range("E3").Value = "  +1.82%  "
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "213.70"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +1.23%  "
$ws.Range("E6").Value = "  +0.62%  "
$ws.Range("E7").Value = "  -0.12%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "23.47"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +2.98%  "
$ws.Range("E9").Value = "  +1.21%  "
$ws.Range("E10").Value = "  +0.38%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0873"
$ws.Range("D11").ClearFormats()
$ws.Range("D12").Value = "1.881.60"
$ws.Range("E12").Value = "  +1.85%  "
$ws.Range("D13").Value = "1.641.88"
$ws.Range("E13").Value = "  +1.39%  "
$ws.Range("E14").Value = "  +1.28%  "
$ws.Range("E15").Value = "  +2.81%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "65.61"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +0.70%  "
$ws.Range("D17").Value = "28.003.83"
$ws.Range("E17").Value = "  +1.95%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "233.15"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +0.97%  "
$ws.Range("E19").Value = "  +2.45%  "
$ws.Range("E22").Value = "  +4.88%  "
$ws.Range("E23").Value = "  +2.87%  "
$ws.Range("E24").Value = "  +4.72%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "152.61"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +1.48%  "
$ws.Range("E26").Value = "  +1.28%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "15.77"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +1.27%  "
$ws.Range("E28").Value = "  +0.37%  "
$ws.Range("E29").Value = "  -0.07%  "
$ws.Range("E30").Value = "  +1.47%  "
$ws.Range("E31").Value = "  +0.05%  "
$ws.Range("E32").Value = "  +2.86%  "
$ws.Range("D33").Value = "1.449.25"
$ws.Range("E33").Value = "  -0.13%  "
$ws.Range("E34").Value = "  +0.85%  "
$ws.Range("E35").Value = "  +1.85%  "
$ws.Range("E36").Value = "  -0.40%  "
$ws.Range("E38").Value = "  +0.91%  "
$ws.Range("E39").Value = "  +0.50%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.921"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -1.97%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "69.49"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +2.54%  "
$ws.Range("E42").Value = "  +3.08%  "
$ws.Range("E43").Value = "  -0.15%  "
$ws.Range("E44").Value = "  -0.43%  "
$ws.Range("E45").Value = "  +0.97%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "5.39"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -0.32%  "
$ws.Range("E47").Value = "  +4.78%  "
$ws.Range("D48").Value = "1.789.81"
$ws.Range("E48").Value = "  +1.56%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "89.05"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +2.90%  "
$ws.Range("E50").Value = "  -1.14%  "
$ws.Range("E51").Value = "  +0.58%  "
